# expression-showcase.xlsx edit script
# Implements:
#  - New "sms" command group (single command: sendText(phones,text)) inserted in column Q
#  - New "ws.async" command group (async variants of ws http verbs) inserted in column Y
#  - Existing command-group columns on the '#system' sheet shift right to make room:
#       Q(sound)->R, R(ssh)->S, S(step)->T, T(web)->U, U(webalert)->V,
#       V(webcookie)->W, W(ws)->X, X(xml)->Z
#  - 'target' column (A) gains two more category names (sms, ws.async) in sorted order
#  - defined names updated to reflect the new ranges; two new defined names added

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Shift data right on the system sheet, reading source values fully into
#    memory first so that overlapping source/destination ranges do not clobber
#    each other.
# ---------------------------------------------------------------------------

# old column X ("xml", rows 1-11) moves two columns right, to Z
$xmlVals = $ws.Range("X1:X11").Value2
$ws.Range("Z1:Z11").Value2 = $xmlVals

# old columns Q..W ("sound","ssh","step","web","webalert","webcookie","ws"; rows 1-108,
# "web" being the tallest list) each move one column right, to R..X
$shiftVals = $ws.Range("Q1:W108").Value2
$ws.Range("R1:X108").Value2 = $shiftVals

# ---------------------------------------------------------------------------
# 2) Clear the old "Q" column (now superseded by the new "sms" column content)
#    beyond what the new sms list occupies.
# ---------------------------------------------------------------------------
$ws.Range("Q1:Q108").ClearContents()

# ---------------------------------------------------------------------------
# 3) Write new "sms" column (Q) and new "ws.async" column (Y)
# ---------------------------------------------------------------------------
$ws.Range("Q1").Value2 = "sms"
$ws.Range("Q2").Value2 = "sendText(phones,text)"

$ws.Range("Y1").Value2 = "ws.async"
$ws.Range("Y2").Value2 = "download(url,queryString,saveTo)"
$ws.Range("Y3").Value2 = "get(url,queryString,output)"
$ws.Range("Y4").Value2 = "head(url,output)"
$ws.Range("Y5").Value2 = "patch(url,body,output)"
$ws.Range("Y6").Value2 = "post(url,body,output)"
$ws.Range("Y7").Value2 = "put(url,body,output)"

# ---------------------------------------------------------------------------
# 4) Extend the "target" list (column A) with the two new category names,
#    keeping the whole list alphabetically sorted.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value2 = "sms"
$ws.Range("A18").Value2 = "sound"
$ws.Range("A19").Value2 = "ssh"
$ws.Range("A20").Value2 = "step"
$ws.Range("A21").Value2 = "web"
$ws.Range("A22").Value2 = "webalert"
$ws.Range("A23").Value2 = "webcookie"
$ws.Range("A24").Value2 = "ws"
$ws.Range("A25").Value2 = "ws.async"
$ws.Range("A26").Value2 = "xml"

# ---------------------------------------------------------------------------
# 5) Update defined names to the new ranges, and add the two new ones.
# ---------------------------------------------------------------------------
$wb.Names.Item("ssh").RefersTo = "='#system'!`$S`$2:`$S`$9"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$26"
$wb.Names.Item("web").RefersTo = "='#system'!`$U`$2:`$U`$108"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$V`$2:`$V`$6"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$X`$2:`$X`$16"
$wb.Names.Item("xml").RefersTo = "='#system'!`$Z`$2:`$Z`$11"
$wb.Names.Item("step").RefersTo = "='#system'!`$T`$2:`$T`$4"
$wb.Names.Item("sound").RefersTo = "='#system'!`$R`$2:`$R`$5"

$wb.Names.Add("sms", "='#system'!`$Q`$2:`$Q`$2")
$wb.Names.Add("ws.async", "='#system'!`$Y`$2:`$Y`$7")
